$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the mapping table from columns C:D to columns A:B -----------------
# Copy the label/value pairs (rows 1-13) from C:D into A:B.
$vals = $ws.Range("C1:D13").Value()
$ws.Range("A1:B13").Value = $vals

# Re-apply the formatting the cells had in their old location:
#  - header row (row 1) used the bold style
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Size = 8
$ws.Range("A1:B1").Font.Name = "Arial"

#  - data rows (2-13), column A used the plain "label" style; column B used
#    the default (unstyled) formatting, so it is left untouched.
$ws.Range("A2:A13").Font.Bold = $false
$ws.Range("A2:A13").Font.Size = 8
$ws.Range("A2:A13").Font.Name = "Arial"

# The old C:D columns are now empty - clear their contents & formatting.
$ws.Range("C1:D13").Clear()

# --- Drop the "YoA" year column (F) from the mapping UI ---------------------
$ws.Range("F1").Value = ""
$ws.Range("F2").Value = ""

# --- Resize the new columns to fit their (now relocated) content ------------
$ws.Columns("A").ColumnWidth = 13.6
$ws.Columns("B").ColumnWidth = 34.3

# --- Update selection to match the saved view --------------------------------
$ws.Range("C17").Select() | Out-Null
